$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as literal text (e.g. "0.660",
# "0.0000300", "1.00") in the source data, so we temporarily force
# text format before assigning them, then restore the default cell
# style so the cell keeps behaving like a normal/unstyled cell.
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '70.138.19'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.90%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.584.28'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.09%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '576.07'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -2.61%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '190.58'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.71%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.633'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.62%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '3.582.13'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  -3.21%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.660'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -0.25%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '56.61'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -2.79%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000300'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +2.76%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '9.80'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +1.04%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '4.151.00'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +0.37%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '20.14'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +4.55%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '3.572.67'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.26%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '70.020.96'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.81%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '12.52'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('E21').Value = '  -0.47%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '19.80'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +15.69%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '472.70'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -6.05%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '5.12'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -7.46%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '4.34'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.23%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '88.68'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -2.76%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '3.06'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +0.12%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '11.15'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.11%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '9.29'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '7.80'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +3.87%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '32.10'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  +5.10%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '12.12'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.18%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '66.26'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +1.35%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '590.68'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -4.55%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '39.65'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +4.44%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0807'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.30%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.401'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '3.55'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '2.93'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +7.92%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '3.231.06'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -2.99%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E44').Value = '  +8.21%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '3.14'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +2.34%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.0447'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +1.00%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '9.64'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +5.86%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '3.34'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +0.98%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.137'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.47%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.997'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -0.30%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '3.17'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -2.12%  '

Write-Host "Applied 100 cell updates"
